$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 of data
$ws.Range("B3").Value = "dummy"

$ws.Range("C3").NumberFormat = "mm-dd-yy"
$ws.Range("C3").Value = Get-Date -Year 2016 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Range("D3").Value = 75519
$ws.Range("E3").Value = "gut"
$ws.Range("F3").Value = "dummy original"
$ws.Range("H3").Value = "dummy original"

# Update selection to match the post-edit state
$ws.Range("B4").Select()
